$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Widen column C to fit the new, longer topic text ---
$ws.Columns.Item(3).ColumnWidth = 67.17

# --- 2. Fill in Topic (column C) values for existing rows that were blank ---
$ws.Range("C15").Value  = "Saturday Holiday"
$ws.Range("C22").Value  = "Saturday Holiday"
$ws.Range("C24").Value  = "Bhatta Sir"
$ws.Range("C25").Value  = "Bhatta Sir"
$ws.Range("C26").Value  = "Bhatta Sir"
$ws.Range("C27").Value  = "Bhatta Sir"
$ws.Range("C28").Value  = "Bhatta Sir"
$ws.Range("C29").Value  = "Saturday Holiday"
$ws.Range("C30").Value  = "Bhatta Sir"
$ws.Range("C31").Value  = "Bhatta Sir"
$ws.Range("C32").Value  = "Bhatta Sir"
$ws.Range("C33").Value  = "Bhatta Sir"
$ws.Range("C34").Value  = "Bhatta Sir"
$ws.Range("C35").Value  = "Bhatta Sir"
$ws.Range("C36").Value  = "Saturday Holiday"
$ws.Range("C37").Value  = "Bhatta Sir"
$ws.Range("C38").Value  = "Bhatta Sir"
$ws.Range("C39").Value  = "Bhatta Sir"
$ws.Range("C40").Value  = "Bhatta Sir"
$ws.Range("C41").Value  = "Bhatta Sir"
$ws.Range("C50").Value  = "Saturday Holiday"
$ws.Range("C55").Value  = "Sports Week"
$ws.Range("C56").Value  = "Sports Week"
$ws.Range("C57").Value  = "Saturday Holiday"
$ws.Range("C58").Value  = "Sports Week"
$ws.Range("C59").Value  = "Ghode Jatra  Holiday"

# --- 3. Row 59 shrinks back to the regular row height ---
$ws.Rows.Item(59).RowHeight = 19.5

# --- 4. Append 19 new rows (60-78) continuing the calendar into April 2024 ---
$ws.Range("A59:C59").Copy() | Out-Null
$ws.Range("A60:C78").PasteSpecial(-4122) | Out-Null

for ($r = 60; $r -le 78; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
    $ws.Cells.Item($r, 2).Value = 45391 + ($r - 60)
    $ws.Rows.Item($r).RowHeight = 18.75
}

# Row 60 carries the last entry text, rows 61-78 stay blank in column C
$ws.Range("C60").Value = "Multimedia elements HTML, audio, video, YouTube and Border Properties"
